# This script applies a row-content rotation to the "Artfynd" sheet.
# Each row keeps its position-dependent columns (P, S, T, U, V, W, Z, AB, C,
# I, AD, AE, AG, AT, AW, AX, AY) but the record-specific columns (A, B, D, E,
# F, G, H, Q, R, Y, AA, AF) rotate between rows as follows:
#   row2  <- row4,  row4  <- row5,  row5  <- row2    (cycle A)
#   row6  <- row9,  row9  <- row11, row11 <- row6    (cycle B)
#   row7  <- row10, row10 <- row7                    (cycle C)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns Y and AA hold plain text that LOOKS like a date ("2023-08-29").
# Reading such a cell's .Value2 and writing it straight back in makes Excel's
# COM layer auto-convert it to a date serial number, so those two columns are
# written back with the cell pre-formatted as Text, then the format is reset
# back to the default "Normal" style so no stray number-format is left behind.
$dateLikeCols = @("Y", "AA")
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Y", "AA")

function Read-RowSnapshot($row) {
    $snap = @{}
    foreach ($c in $cols) {
        $addr = "$c$row"
        $snap[$c] = $ws.Range($addr).Value2
    }
    return $snap
}

function Write-RowSnapshot($row, $snap) {
    foreach ($c in $cols) {
        $addr = "$c$row"
        if ($dateLikeCols -contains $c) {
            $ws.Range($addr).NumberFormat = "@"
            $ws.Range($addr).Value = $snap[$c]
            $ws.Range($addr).Style = "Normal"
        } else {
            $ws.Range($addr).Value = $snap[$c]
        }
    }
}

# --- snapshot every source row BEFORE any writes happen ---
$snap2  = Read-RowSnapshot 2
$snap4  = Read-RowSnapshot 4
$snap5  = Read-RowSnapshot 5
$snap6  = Read-RowSnapshot 6
$snap7  = Read-RowSnapshot 7
$snap9  = Read-RowSnapshot 9
$snap10 = Read-RowSnapshot 10
$snap11 = Read-RowSnapshot 11

# --- apply the rotations ---
# cycle A: 2 <- 4 <- 5 <- 2
Write-RowSnapshot 2 $snap4
Write-RowSnapshot 4 $snap5
Write-RowSnapshot 5 $snap2

# cycle B: 6 <- 9 <- 11 <- 6
Write-RowSnapshot 6 $snap9
Write-RowSnapshot 9 $snap11
Write-RowSnapshot 11 $snap6

# cycle C: 7 <- 10 <- 7
Write-RowSnapshot 7 $snap10
Write-RowSnapshot 10 $snap7

# --- AF column (Bestamningsmetod) travels with the record too. ---
# Before the edit, AF10 and AF11 hold a present-but-empty value while AF7 and
# AF9 are blank; after the edit AF7 (<- old AF10) and AF9 (<- old AF11) gain
# that present-but-empty value, and AF10 / AF11 go back to being blank.

# AF7 / AF9 must become present-but-empty text cells (they were blank before).
$ws.Range("AF7").Value = "'"
$ws.Range("AF7").Style = "Normal"
$ws.Range("AF9").Value = "'"
$ws.Range("AF9").Style = "Normal"

# AF10 / AF11 were present-but-empty text cells; they must become blank again.
$ws.Range("AF10").ClearContents()
$ws.Range("AF11").ClearContents()
